$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$origStyle = $ws.Range("D2").Style
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '53.274.07'
$ws.Range("D2").Style = $origStyle
$ws.Range("E2").Value = '  -5.06%  '

$origStyle = $ws.Range("D3").Style
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.180.93'
$ws.Range("D3").Style = $origStyle
$ws.Range("E3").Value = '  -7.81%  '

$origStyle = $ws.Range("D4").Style
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = $origStyle
$ws.Range("E4").Value = '  -0.15%  '

$origStyle = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '480.14'
$ws.Range("D5").Style = $origStyle
$ws.Range("E5").Value = '  -4.22%  '

$origStyle = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '123.39'
$ws.Range("D6").Style = $origStyle
$ws.Range("E6").Value = '  -4.43%  '

$origStyle = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.997'
$ws.Range("D7").Style = $origStyle
$ws.Range("E7").Value = '  -0.26%  '

$origStyle = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.515'
$ws.Range("D8").Style = $origStyle
$ws.Range("E8").Value = '  -5.35%  '

$origStyle = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.195.03'
$ws.Range("D9").Style = $origStyle
$ws.Range("E9").Value = '  -7.16%  '

$origStyle = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0908'
$ws.Range("D10").Style = $origStyle
$ws.Range("E10").Value = '  -7.42%  '

$ws.Range("E11").Value = '  -2.00%  '

$origStyle = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.56'
$ws.Range("D12").Style = $origStyle
$ws.Range("E12").Value = '  -5.56%  '

$origStyle = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.311'
$ws.Range("D13").Style = $origStyle
$ws.Range("E13").Value = '  -3.52%  '

$origStyle = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.569.02'
$ws.Range("D14").Style = $origStyle
$ws.Range("E14").Value = '  -7.85%  '

$origStyle = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '20.85'
$ws.Range("D15").Style = $origStyle
$ws.Range("E15").Value = '  -2.61%  '

$origStyle = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '53.182.32'
$ws.Range("D16").Style = $origStyle
$ws.Range("E16").Value = '  -5.19%  '

$ws.Range("E17").Value = '  -4.22%  '

$origStyle = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.178.23'
$ws.Range("D18").Style = $origStyle
$ws.Range("E18").Value = '  -5.07%  '

$origStyle = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '9.48'
$ws.Range("D19").Style = $origStyle
$ws.Range("E19").Value = '  -5.20%  '

$origStyle = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '3.91'
$ws.Range("D20").Style = $origStyle
$ws.Range("E20").Value = '  -3.15%  '

$origStyle = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '290.84'
$ws.Range("D21").Style = $origStyle
$ws.Range("E21").Value = '  -5.08%  '

$origStyle = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.999'
$ws.Range("D23").Style = $origStyle
$ws.Range("E23").Value = '  -0.02%  '

$origStyle = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '62.14'
$ws.Range("D24").Style = $origStyle
$ws.Range("E24").Value = '  -5.31%  '

$origStyle = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.997'
$ws.Range("D25").Style = $origStyle
$ws.Range("E25").Value = '  +0.03%  '

$ws.Range("E26").Value = '  -2.43%  '

$origStyle = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.282.97'
$ws.Range("D27").Style = $origStyle
$ws.Range("E27").Value = '  -7.81%  '

$ws.Range("E28").Value = '  -2.44%  '

$ws.Range("E29").Value = '  -3.75%  '

$origStyle = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '165.39'
$ws.Range("D30").Style = $origStyle
$ws.Range("E30").Value = '  -3.53%  '

$ws.Range("E32").Value = '  -4.69%  '

$ws.Range("E33").Value = '  -0.33%  '

$origStyle = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0₃0654'
$ws.Range("D34").Style = $origStyle
$ws.Range("E34").Value = '  -7.87%  '

$origStyle = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.61'
$ws.Range("D35").Style = $origStyle
$ws.Range("E35").Value = '  -2.28%  '

$ws.Range("E36").Value = '  -3.52%  '

$origStyle = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '17.15'
$ws.Range("D37").Style = $origStyle
$ws.Range("E37").Value = '  -2.62%  '

$origStyle = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.14'
$ws.Range("D38").Style = $origStyle
$ws.Range("E38").Value = '  -2.83%  '

$origStyle = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.806'
$ws.Range("D39").Style = $origStyle
$ws.Range("E39").Value = '  +1.98%  '

$origStyle = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '35.59'
$ws.Range("D40").Style = $origStyle
$ws.Range("E40").Value = '  -1.34%  '

$ws.Range("E41").Value = '  -6.23%  '

$ws.Range("E42").Value = '  -2.21%  '

$origStyle = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.35'
$ws.Range("D43").Style = $origStyle
$ws.Range("E43").Value = '  -2.47%  '

$origStyle = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.22'
$ws.Range("D44").Style = $origStyle
$ws.Range("E44").Value = '  -3.72%  '

$origStyle = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '122.39'
$ws.Range("D45").Style = $origStyle
$ws.Range("E45").Value = '  -5.10%  '

$origStyle = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '4.68'
$ws.Range("D46").Style = $origStyle
$ws.Range("E46").Value = '  -0.35%  '

$origStyle = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0871'
$ws.Range("D47").Style = $origStyle
$ws.Range("E47").Value = '  -3.41%  '

$ws.Range("E48").Value = '  -6.27%  '

$origStyle = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0465'
$ws.Range("D49").Style = $origStyle
$ws.Range("E49").Value = '  -3.25%  '

$origStyle = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '226.82'
$ws.Range("D50").Style = $origStyle
$ws.Range("E50").Value = '  -5.09%  '

$ws.Range("E51").Value = '  -4.30%  '
